# Apply the changes described by the diff:
# 1. Clear the stray empty cell B6 on the "ODI Batting" sheet.
# 2. Add a new worksheet "ODI Batting Extra" at the end of the workbook
#    with MATCH_CODE / BATTING_POSITION / NUM_4 / NUM_6 /
#    PERCENT_RUNS_OF_TOTAL / MAN_OF_MATCH columns.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. ODI Batting: clear the empty inline-string cell left in B6
# ---------------------------------------------------------------------
$odiBatting = $wb.Worksheets.Item("ODI Batting")
$odiBatting.Range("B6").ClearContents()

# ---------------------------------------------------------------------
# 2. Add the new "ODI Batting Extra" worksheet after the last sheet
# ---------------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$newSheet.Name = "ODI Batting Extra"

# Copy the header style (bold, centered, bordered) from an existing
# sheet's header row so the new header reuses the same cell style.
$srcHeader = $odiBatting.Range("A1:F1")
$srcHeader.Copy($newSheet.Range("A1:F1"))

$newSheet.Range("A1").Value = "MATCH_CODE"
$newSheet.Range("B1").Value = "BATTING_POSITION"
$newSheet.Range("C1").Value = "NUM_4"
$newSheet.Range("D1").Value = "NUM_6"
$newSheet.Range("E1").Value = "PERCENT_RUNS_OF_TOTAL"
$newSheet.Range("F1").Value = "MAN_OF_MATCH"

# Force the data block to be treated as text so values such as "4686",
# "0", and "1.85%" are stored as strings rather than numbers/percentages.
$newSheet.Range("A2:F6").NumberFormat = "@"

$newSheet.Range("A2").Value = "4686"
$newSheet.Range("F2").Value = "NO"

$newSheet.Range("A3").Value = "4692"
$newSheet.Range("F3").Value = "NO"

$newSheet.Range("A4").Value = "4695"
$newSheet.Range("C4").Value = "0"
$newSheet.Range("D4").Value = "0"
$newSheet.Range("E4").Value = "1.85%"
$newSheet.Range("F4").Value = "NO"

$newSheet.Range("A5").Value = "4735"
$newSheet.Range("C5").Value = "0"
$newSheet.Range("D5").Value = "0"
$newSheet.Range("E5").Value = "2.19%"
$newSheet.Range("F5").Value = "YES"

$newSheet.Range("A6").Value = "4745"
$newSheet.Range("F6").Value = "NO"

# B4 and B5 are actual numeric batting positions (9 and 8), so set them
# after the text formatting was applied to the rest of the block, and
# clear the text format from just these two cells so they serialize as
# numbers.
$newSheet.Range("B4").NumberFormat = "General"
$newSheet.Range("B4").Value = 9
$newSheet.Range("B5").NumberFormat = "General"
$newSheet.Range("B5").Value = 8
